$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6539
$ws.Range("C27").Value = 1019
$ws.Range("D27").Value = 6099083
$ws.Range("E27").Value = 932.7241168374369
$ws.Range("F27").Value = 9.899159663865543
$ws.Range("G27").Value = 7.602956705385422
$ws.Range("H27").Value = 25.0963848477707
